$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("G6").Value = 1.5
$ws.Range("I6").Value = 6.5
$ws.Range("J6").Value = 2.1
$ws.Range("L6").Value = 6.5
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("Y6").Value = 8.5
$ws.Range("Z6").Value = 10
$ws.Range("AB6").Value = 29
$ws.Range("AK6").Value = 67
$ws.Range("AN6").Value = 3.4
$ws.Range("AS6").Value = 151
$ws.Range("AW6").Value = 7.5
$ws.Range("AX6").Value = 34
$ws.Range("BA6").Value = 151
$ws.Range("G8").Value = 1.4
$ws.Range("H8").Value = 4.33
$ws.Range("I8").Value = 8.5
$ws.Range("J8").Value = 1.91
$ws.Range("L8").Value = 7.5
$ws.Range("Z8").Value = 9
$ws.Range("AD8").Value = 8.5
$ws.Range("AE8").Value = 21
$ws.Range("AQ8").Value = 19
$ws.Range("BA8").Value = 201
$ws.Range("J9").Value = 1.95
$ws.Range("L9").Value = 8
$ws.Range("M9").Value = 1.06
$ws.Range("N9").Value = 10
$ws.Range("Q9").Value = 2.03
$ws.Range("R9").Value = 1.83
$ws.Range("U9").Value = 2.25
$ws.Range("V9").Value = 1.57
$ws.Range("W9").Value = 5.5
$ws.Range("AC9").Value = 9
$ws.Range("AE9").Value = 23
$ws.Range("AJ9").Value = 26
$ws.Range("BA9").Value = 251
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
$ws.Range("O11").Value = 1.33
$ws.Range("P11").Value = 3.25
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 4.5
$ws.Range("I14").Value = 1.44
$ws.Range("K14").Value = 2.3
$ws.Range("L14").Value = 2
$ws.Range("Q14").Value = 1.95
$ws.Range("R14").Value = 1.9
$ws.Range("AH14").Value = 6.5
$ws.Range("AJ14").Value = 9
$ws.Range("AM14").Value = 34
$ws.Range("AR14").Value = 201
$ws.Range("AX14").Value = 7
$ws.Range("AZ14").Value = 21
